# Apply updated crypto price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.908.56"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "1.636.41"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'215.41"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").Value = "'19.54"
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.863.53"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.27"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").Value = "1.639.93"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "25.949.48"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "'192.90"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("D22").Value = "'9.90"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("E24").Value = "  +4.07%  "
$ws.Range("D25").Value = "'1.79"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'143.15"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").Value = "'15.54"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -4.27%  "
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("D36").Value = "'0.899"
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("D37").Value = "1.133.14"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").Value = "'2.45"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").Value = "1.773.23"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").Value = "'56.53"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  -1.54%  "
